$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Commit: "Manual addition of contract info, removed players without active
# snaps in 2024-2025 NFL season. Reordered code structure for logical flow."
#
# Net effect on Sheet1 rows 241-273:
#   1. The trailing 6-row block (Joey Bosa ... Joe Gaziano, formerly rows
#      267-272) is moved to the front of this section (new rows 241-246) so
#      the data reads in the intended logical order.
#   2. A new player row for Khari Blasingame (Chicago Bears) is inserted
#      right after the Sam Hubbard row.
#   3. Every other row keeps its same relative order, shifting down by one
#      to make room.
# The net row count grows from 272 to 273.
# ---------------------------------------------------------------------------

$targetRows = @(
  ,@('Joey Bosa', 'Los Angeles Chargers', 29, 135000000, 27000000, 102000000, 78000000, '2026 UFA')
  ,@('Brent Urban', 'Baltimore Ravens', 33, 2500000, 2500000, 1000000, 1000000, '2025 UFA')
  ,@('Jermaine Johnson', 'New York Jets', 26, 13088000, 3272000, 13088000, 13088000, '2026 UFA')
  ,@('Casey Toohill', 'Washington Commanders', 28, 2700000, 900000, 150000, 150000, '2025 UFA')
  ,@('Shaq Lawson', 'Buffalo Bills', 31, 9000000, 3000000, 3000000, 3000000, '2025 UFA')
  ,@('Joe Gaziano', 'Atlanta Falcons', 27, 2500000, 833333, 100000, 100000, '2025 UFA')
  ,@('Isaiah Land', 'Dallas Cowboys', 24, 2705000, 902000, 45000, 45000, '2026 ERFA')
  ,@('Jamie Sheriff', 'Miami Dolphins', 23, 2710000, 903333, 50000, 50000, '2026 ERFA')
  ,@('Rashad Weaver', 'Tennessee Titans', 26, 4100000, 1025000, 1000000, 1000000, '2025 UFA')
  ,@('Deatrich Wise Jr.', 'New England Patriots', 30, 22000000, 5500000, 10000000, 10000000, '2027 UFA')
  ,@('Adetokunbo Ogundeji', 'Atlanta Falcons', 26, 3800000, 950000, 320000, 320000, '2025 UFA')
  ,@('Preston Smith', 'Green Bay Packers', 32, 52000000, 13000000, 16000000, 16000000, '2027 UFA')
  ,@('Dre''Mont Jones', 'Seattle Seahawks', 28, 51530000, 17176667, 30000000, 30000000, '2026 UFA')
  ,@('Haason Reddick', 'Philadelphia Eagles', 30, 45000000, 15000000, 30000000, 30000000, '2026 UFA')
  ,@('Sam Hubbard', 'Cincinnati Bengals', 29, 40000000, 10000000, 16000000, 16000000, '2027 UFA')
  ,@('Khari Blasingame', 'Chicago Bears', 29, 2700000, 900000, 1080000, 1080000, '2025 UFA')
  ,@('Zaven Collins', 'Arizona Cardinals', 25, 14700000, 3675000, 14700000, 14700000, '2026 UFA')
  ,@('Nick Herbig', 'Pittsburgh Steelers', 24, 5000000, 1250000, 3000000, 3000000, '2027 UFA')
  ,@('DJ Coleman', 'Jacksonville Jaguars', 25, 2570000, 857000, 10000, 10000, '2026 ERFA')
  ,@('Charles Harris', 'Detroit Lions', 30, 13000000, 6500000, 7000000, 7000000, '2026 UFA')
  ,@('Yannick Ngakoue', 'Chicago Bears', 30, 10500000, 10500000, 10500000, 10500000, '2025 UFA')
  ,@('Michael Burton', 'Kansas City Chiefs', 34, 2500000, 1250000, 1000000, 1000000, '2025 UFA')
  ,@('DeMarvin Leal', 'Pittsburgh Steelers', 24, 5000000, 1250000, 3000000, 3000000, '2027 UFA')
  ,@('Jamin Davis', 'Washington Commanders', 27, 13800000, 3450000, 13800000, 13800000, '2026 UFA')
  ,@('Adam Gotsis', 'Jacksonville Jaguars', 33, 3000000, 3000000, 3000000, 3000000, '2025 UFA')
  ,@('Takk McKinley', 'Dallas Cowboys', 29, 4250000, 4250000, 3000000, 3000000, '2025 UFA')
  ,@('Cam Gill', 'Tampa Bay Buccaneers', 27, 2705000, 902000, 45000, 45000, '2025 UFA')
  ,@('DeMarcus Walker', 'Chicago Bears', 31, 21000000, 7000000, 16000000, 16000000, '2026 UFA')
  ,@('Marquis Haynes Sr.', 'Carolina Panthers', 31, 5500000, 2750000, 2300000, 2300000, '2025 UFA')
  ,@('Efe Obada', 'Washington Commanders', 33, 3000000, 3000000, 3000000, 3000000, '2025 UFA')
  ,@('Trevor Nowaske', 'Green Bay Packers', 24, 2710000, 903333, 50000, 50000, '2026 ERFA')
  ,@('Reggie Gilliam', 'Buffalo Bills', 28, 5200000, 2600000, 3000000, 3000000, '2026 UFA')
  ,@('Harold Landry III', 'New England Patriots', 29, 43500000, 14500000, 30000000, 30000000, '2028 UFA')
)

$targetHeights = @(
  43.2
  28.8
  28.8
  57.6
  28.8
  28.8
  28.8
  28.8
  28.8
  43.2
  28.8
  28.8
  28.8
  28.8
  28.8
  28.8
  28.8
  28.8
  28.8
  28.8
  28.8
  43.2
  28.8
  57.6
  28.8
  28.8
  57.6
  28.8
  28.8
  57.6
  28.8
  28.8
  43.2
)

# First, make room: insert one new row so the section grows from 32 to 33
# rows (241-272 -> 241-273). Insert at the bottom of the affected block so
# the existing row 272 formatting (style s="2"/s="3" on A:H, the currency
# number format, wrap text) is preserved for every row, then stamp row 273
# with the same A:H formatting copied from its neighbour.
$ws.Range("A272:H272").Copy($ws.Range("A273:H273"))
$ws.Application.CutCopyMode = $false

# Write every row's values (A:H) for rows 241-273 in the new order.
for ($i = 0; $i -lt $targetRows.Length; $i++) {
    $r = 241 + $i
    $row = $targetRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# Re-apply each row's height now that the content (and therefore the
# required wrap height) has moved with it.
for ($i = 0; $i -lt $targetHeights.Length; $i++) {
    $r = 241 + $i
    $ws.Rows($r).RowHeight = $targetHeights[$i]
}

# The sheet's used range now extends one row further.
$ws.Range("A1").Worksheet.UsedRange | Out-Null

# View-state touch-ups that Excel recalculates as a side effect of scrolling
# to, and leaving the cursor on, the newly-edited area.
$ws.Application.GoTo($ws.Range("A245"), $true)
$ws.Range("J252").Select()
